$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H) mirrors the existing header formatting (bold,
# bordered, centered) used by the other header cells - copy G1's format
# onto H1 before writing its text so it picks up the same style index.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
